# Auto-generated Excel COM-interop script
# Applies numeric value updates to the Typhon_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as produced by the scheduled market-data refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H6").Value = 516.6667
$ws.Range("J6").Value = 750
$ws.Range("L6").Value = 2250
$ws.Range("N6").Value = -2474
$ws.Range("H18").Value = 159.96552
$ws.Range("I18").Value = 159.96552
$ws.Range("K18").Value = 159.96552
$ws.Range("M18").Value = 124.03448
$ws.Range("H38").Value = 698.8
$ws.Range("I38").Value = 62.666668
$ws.Range("J38").Value = 971.4286
$ws.Range("K38").Value = 188.000004
$ws.Range("L38").Value = 2914.2858
$ws.Range("M38").Value = 183.999996
$ws.Range("N38").Value = -3658.2858
$ws.Range("H39").Value = 255.27272
$ws.Range("I39").Value = 214
$ws.Range("J39").Value = 327.5
$ws.Range("K39").Value = 642
$ws.Range("L39").Value = 982.5
$ws.Range("M39").Value = -346
$ws.Range("N39").Value = -1574.5
$ws.Range("H62").Value = 2651.8333
$ws.Range("I62").Value = 1920.8572
$ws.Range("J62").Value = 3675.2
$ws.Range("K62").Value = 1920.8572
$ws.Range("L62").Value = 3675.2
$ws.Range("M62").Value = -1296.8572
$ws.Range("N62").Value = -4923.2
$ws.Range("H65").Value = 2651.8333
$ws.Range("I65").Value = 1920.8572
$ws.Range("J65").Value = 3675.2
$ws.Range("K65").Value = 9604.286
$ws.Range("L65").Value = 18376
$ws.Range("M65").Value = -6484.286
$ws.Range("N65").Value = -24616
$ws.Range("H69").Value = 1490.7142
$ws.Range("J69").Value = 1466.9117
$ws.Range("L69").Value = 4400.7351
$ws.Range("N69").Value = -6148.7351
$ws.Range("H72").Value = 1490.7142
$ws.Range("J72").Value = 1466.9117
$ws.Range("L72").Value = 13202.2053
$ws.Range("N72").Value = -21938.2053
$ws.Range("H100").Value = 2924
$ws.Range("I100").Value = 603.3333
$ws.Range("J100").Value = 4084.3333
$ws.Range("K100").Value = 603.3333
$ws.Range("L100").Value = 4084.3333
$ws.Range("M100").Value = -62.33330000000001
$ws.Range("N100").Value = -5166.3333
$ws.Range("H106").Value = 2002.4642
$ws.Range("I106").Value = 1591.762
$ws.Range("K106").Value = 1591.762
$ws.Range("M106").Value = -960.7619999999999
$ws.Range("H107").Value = 520.9474
$ws.Range("I107").Value = 482.5
$ws.Range("J107").Value = 586.8570999999999
$ws.Range("K107").Value = 482.5
$ws.Range("L107").Value = 586.8570999999999
$ws.Range("M107").Value = 1437.5
$ws.Range("N107").Value = -4426.8571
$ws.Range("H138").Value = 1678.52
$ws.Range("I138").Value = 646.8
$ws.Range("J138").Value = 2120.6858
$ws.Range("K138").Value = 1940.4
$ws.Range("L138").Value = 6362.057400000001
$ws.Range("M138").Value = 3199.6
$ws.Range("N138").Value = -16642.0574
$ws.Range("H141").Value = 2052.0476
$ws.Range("I141").Value = 1587.5294
$ws.Range("J141").Value = 4026.25
$ws.Range("K141").Value = 4762.5882
$ws.Range("L141").Value = 12078.75
$ws.Range("M141").Value = 417.4117999999999
$ws.Range("N141").Value = -22438.75

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 2493.4546
$ws.Range("I32").Value = 2071.4263
$ws.Range("K32").Value = 2071.4263
$ws.Range("M32").Value = -1784.4263
$ws.Range("H74").Value = 66667828
$ws.Range("I74").Value = 100000700
$ws.Range("J74").Value = 2079.6
$ws.Range("K74").Value = 100000700
$ws.Range("L74").Value = 2079.6
$ws.Range("M74").Value = -99999826
$ws.Range("N74").Value = -3827.6
$ws.Range("H77").Value = 66667828
$ws.Range("I77").Value = 100000700
$ws.Range("J77").Value = 2079.6
$ws.Range("K77").Value = 500003500
$ws.Range("L77").Value = 10398
$ws.Range("M77").Value = -499999132
$ws.Range("N77").Value = -19134
$ws.Range("H122").Value = 2007.8334
$ws.Range("I122").Value = 1909.4
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 5728.200000000001
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -3278.200000000001
$ws.Range("N122").Value = -12400
$ws.Range("H132").Value = 10628.852
$ws.Range("I132").Value = 1209.4147
$ws.Range("K132").Value = 3628.2441
$ws.Range("M132").Value = -1098.2441

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H134").Value = 3811.6453
$ws.Range("I134").Value = 4301.269
$ws.Range("J134").Value = 1265.6
$ws.Range("K134").Value = 12903.807
$ws.Range("L134").Value = 3796.8
$ws.Range("M134").Value = -10368.807
$ws.Range("N134").Value = -8866.799999999999

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 10632.881
$ws.Range("I31").Value = 19424.611
$ws.Range("J31").Value = 4039.0833
$ws.Range("K31").Value = 19424.611
$ws.Range("L31").Value = 4039.0833
$ws.Range("M31").Value = -19129.611
$ws.Range("N31").Value = -4629.0833
$ws.Range("H34").Value = 10632.881
$ws.Range("I34").Value = 19424.611
$ws.Range("J34").Value = 4039.0833
$ws.Range("K34").Value = 19424.611
$ws.Range("L34").Value = 4039.0833
$ws.Range("M34").Value = -19222.611
$ws.Range("N34").Value = -4443.0833

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H4").Value = 2727354.8
$ws.Range("I4").Value = 93.333336
$ws.Range("J4").Value = 6000068.5
$ws.Range("K4").Value = 280.000008
$ws.Range("L4").Value = 18000205.5
$ws.Range("M4").Value = -168.000008
$ws.Range("N4").Value = -18000429.5
$ws.Range("H6").Value = 158.33333
$ws.Range("I6").Value = 100
$ws.Range("J6").Value = 275
$ws.Range("K6").Value = 300
$ws.Range("L6").Value = 825
$ws.Range("M6").Value = -187
$ws.Range("N6").Value = -1051
$ws.Range("H12").Value = 88.36364
$ws.Range("I12").Value = 49.5
$ws.Range("J12").Value = 97
$ws.Range("K12").Value = 148.5
$ws.Range("L12").Value = 291
$ws.Range("M12").Value = 24.5
$ws.Range("N12").Value = -637
$ws.Range("H131").Value = 787.6900000000001
$ws.Range("J131").Value = 787.6900000000001
$ws.Range("L131").Value = 2363.07
$ws.Range("N131").Value = -12443.07

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H10").Value = 126074.75
$ws.Range("J10").Value = 126074.75
$ws.Range("L10").Value = 126074.75
$ws.Range("N10").Value = -126412.75
$ws.Range("H70").Value = 3686067
$ws.Range("I70").Value = 14576.444
$ws.Range("J70").Value = 7816493.5
$ws.Range("K70").Value = 14576.444
$ws.Range("L70").Value = 7816493.5
$ws.Range("M70").Value = -14306.444
$ws.Range("N70").Value = -7817033.5
$ws.Range("H73").Value = 3686067
$ws.Range("I73").Value = 14576.444
$ws.Range("J73").Value = 7816493.5
$ws.Range("K73").Value = 14576.444
$ws.Range("L73").Value = 7816493.5
$ws.Range("M73").Value = -13640.444
$ws.Range("N73").Value = -7818365.5
$ws.Range("H80").Value = 3435.7693
$ws.Range("I80").Value = 3168.182
$ws.Range("J80").Value = 3632
$ws.Range("K80").Value = 3168.182
$ws.Range("L80").Value = 3632
$ws.Range("M80").Value = -2170.182
$ws.Range("N80").Value = -5628
$ws.Range("H83").Value = 3435.7693
$ws.Range("I83").Value = 3168.182
$ws.Range("J83").Value = 3632
$ws.Range("K83").Value = 15840.91
$ws.Range("L83").Value = 18160
$ws.Range("M83").Value = -10848.91
$ws.Range("N83").Value = -28144
$ws.Range("H132").Value = 34027
$ws.Range("I132").Value = 8765.5
$ws.Range("J132").Value = 54236.2
$ws.Range("K132").Value = 26296.5
$ws.Range("L132").Value = 162708.6
$ws.Range("M132").Value = -23766.5
$ws.Range("N132").Value = -167768.6

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H61").Value = 3824.5715
$ws.Range("I61").Value = 2375
$ws.Range("K61").Value = 2375
$ws.Range("M61").Value = -2173
$ws.Range("H76").Value = 11333
$ws.Range("J76").Value = 11333
$ws.Range("L76").Value = 11333
$ws.Range("N76").Value = -12009
$ws.Range("H79").Value = 11333
$ws.Range("J79").Value = 11333
$ws.Range("L79").Value = 11333
$ws.Range("N79").Value = -13673
$ws.Range("H113").Value = 3824.5715
$ws.Range("I113").Value = 2375
$ws.Range("K113").Value = 2375
$ws.Range("M113").Value = -205

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H132").Value = 1144.1936
$ws.Range("I132").Value = 798.1
$ws.Range("J132").Value = 1773.4546
$ws.Range("K132").Value = 2394.3
$ws.Range("L132").Value = 5320.3638
$ws.Range("M132").Value = 135.6999999999998
$ws.Range("N132").Value = -10380.3638
